$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 37-46 (data no longer produced after perf improvement)
$ws.Range("A37:D46").EntireRow.Delete()

# Updated simulation values for rows 2-36 (columns B, C, D)
$newValues = @{
    2  = @(1.179489608877793, 5.121759991662863, 0.2659999116465858)
    3  = @(2.256852830972132, 2.512543630289937, 0.2822758317278459)
    4  = @(5.890595381288748, 4.54284386572298, 0.6709362115330153)
    5  = @(7.201341856862103, 5.481674990207025, 0.5991265536376147)
    6  = @(22.24114521373448, 5.638927957281405, 0.4474741475953582)
    7  = @(22.3448467568353, 6.69261740908071, 0.5371195989525227)
    8  = @(25.36507346577119, 2.865373845761207, 0.4781371669887253)
    9  = @(36.82177384762808, 4.692078733085003, 0.4321510150696184)
    10 = @(37.18347464075044, 3.82334544200809, 0.5591234724162917)
    11 = @(40.58471787831015, 9.738207570106812, 0.4954519946357452)
    12 = @(41.46396354355726, 2.946521053563037, 0.5109712905263452)
    13 = @(45.12639147527808, 1.622071569012757, 0.3317644782333493)
    14 = @(45.20957483841438, 2.50819376593961, 0.5676103201055435)
    15 = @(47.08168723030688, 1.660466874886036, 0.3523043880356628)
    16 = @(48.53994855760722, 5.517640875896456, 0.5348537277419639)
    17 = @(49.01590151754824, 3.295977289264159, 0.493414690519704)
    18 = @(50.17603911195334, 3.686689180786403, 0.214035410505218)
    19 = @(53.74687108615088, 3.821728631043711, 0.6569225276145009)
    20 = @(54.5193984361466, 2.677302669366846, 0.2952695612447521)
    21 = @(58.18164324253958, 4.656248355095755, 0.311655496227373)
    22 = @(60.79630555341398, 3.817207052540798, 0.6322041413680677)
    23 = @(69.71745092868154, 6.820560017589059, 0.7704641975613831)
    24 = @(70.74070756638815, 7.384142606715507, 0.5529769371399345)
    25 = @(73.39140359159852, 5.292404548805142, 0.5560396686189114)
    26 = @(73.58872520003976, 4.552860150686686, 0.5410710319753527)
    27 = @(77.52095413455471, 2.25718793433356, 0.3963459084395111)
    28 = @(80.27872689282572, 5.856256610519349, 0.6647591598107228)
    29 = @(81.68938598601198, 7.466299566919523, 0.5933720568233973)
    30 = @(83.47866707827706, 6.090542712802796, 0.3560199872119551)
    31 = @(84.58597524997717, 4.297478148273033, 0.4492454992715766)
    32 = @(84.97682422020689, 5.839490860769521, 0.7898600287208435)
    33 = @(89.80862695547808, 3.932820984197781, 0.6399189739808473)
    34 = @(92.39631850671205, 5.368585616669927, 0.2951896464443921)
    35 = @(95.07503045695304, 6.610110400273894, 0.5849567390648386)
    36 = @(99.94446773238099, 3.043545565757124, 0.4845145811264198)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
}
